# BitcoinPriceTracker/Data/prices.xlsx refactor:
#  - "Old Prices" (D) column is shifted: each row's NEW "Old Prices" value
#    becomes the PREVIOUS "Prices" (B) reading for that row.
#  - "Prices" (B) and "Euro" (C) columns get freshly fetched values.
# Row 6 (CoinGecko) stores its numbers as literal text (shared-string,
# trailing NBSP) instead of numeric cells, same as the original workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nbsp = [string][char]0xA0

function Set-TextValue {
    param([string]$addr, [string]$text)
    # Writing a numeric-looking string via .Value auto-converts it to a
    # number, so stage it through a cell formatted as Text ("@") and then
    # copy formatting back off (reusing the existing General-formatted A1
    # cell) so the target cell keeps its original (no) style.
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $ws.Range("A1").Copy() | Out-Null
    $cell.PasteSpecial(-4122) | Out-Null  # xlPasteFormats
}

# Row 2 (Binance): D gets old B, B/C get new readings
$ws.Range("D2").Value = 46080.68
$ws.Range("B2").Value = 46071.03
$ws.Range("C2").Value = 40562.879999999997

# Row 3 (KuCoin)
$ws.Range("D3").Value = 46071.19
$ws.Range("B3").Value = 46066.400000000001
$ws.Range("C3").Value = 40558.81

# Row 4 (Coinbase)
$ws.Range("D4").Value = 46087.23
$ws.Range("B4").Value = 46094.14
$ws.Range("C4").Value = 40583.230000000003

# Row 5 (CMC)
$ws.Range("D5").Value = 46170.09
$ws.Range("B5").Value = 46145.919999999998
$ws.Range("C5").Value = 40628.82

# Row 6 (CoinGecko): text cells (trailing NBSP), plus numeric Euro cell
Set-TextValue "D6" ("46155.04" + $nbsp)
Set-TextValue "B6" ("46131.60" + $nbsp)
$ws.Range("C6").Value = 40616.21

$excel.CutCopyMode = 0
